# This edit reorders the data rows (A:F) on the active sheet so that each
# target row receives the values that originally lived in a different row.
# Rows 5, 6, and 21-26 are left untouched; only rows 3,4,7-20 are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    3  = @(1202, 2, 10, 10, 10, 10)
    4  = @(101, 9, 30, 15, 60, 15)
    7  = @(301, 6, 45, 30, 60, 45)
    8  = @(701, 3, 90, 45, 97, 15)
    9  = @(1201, 2, 10, 10, 10, 10)
    10 = @(201, 9, 30, 15, 45, 30)
    11 = @(801, 3, 67, 65, 52, 45)
    12 = @(501, 9, 52, 30, 75, 45)
    13 = @(401, 9, 48, 67, 75, 45)
    14 = @(601, 9, 60, 67, 60, 42)
    15 = @(1203, 3, 15, 15, 15, 15)
    16 = @(802, 0, 4, 5, 4, 0)
    17 = @(1, 0, 2, 2, 2, 2)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(3, 0, 3, 3, 3, 3)
    20 = @(502, 0, 4, 0, 0, 0)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
